# Updated cryptos list with refreshed price/volume(1h) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.198.03"
$ws.Range("E2").Value = "  -0.93%  "

# Row 3
$ws.Range("D3").Value = "2.431.32"
$ws.Range("E3").Value = "  -1.53%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.25"
$ws.Range("E5").Value = "  +0.14%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "89.74"
$ws.Range("E6").Value = "  -2.48%  "

# Row 7
$ws.Range("E7").Value = "  -2.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("E9").Value = "  -2.98%  "

# Row 10
$ws.Range("E10").Value = "  -1.49%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.15"
$ws.Range("E11").Value = "  -2.31%  "

# Row 12
$ws.Range("E12").Value = "  -1.46%  "

# Row 13
$ws.Range("D13").Value = "2.807.82"
$ws.Range("E13").Value = "  -1.47%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.74"
$ws.Range("E14").Value = "  -1.75%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.73"
$ws.Range("E15").Value = "  +1.37%  "

# Row 16
$ws.Range("D16").Value = "2.417.74"
$ws.Range("E16").Value = "  -2.24%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.776"
$ws.Range("E17").Value = "  -1.74%  "

# Row 18
$ws.Range("D18").Value = "41.117.96"
$ws.Range("E18").Value = "  -1.01%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0928"
$ws.Range("E19").Value = "  -1.87%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.29"
$ws.Range("E20").Value = "  -2.27%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.21"
$ws.Range("E21").Value = "  +1.62%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.10"
$ws.Range("E22").Value = "  -1.57%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.63"
$ws.Range("E23").Value = "  -1.33%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.69"
$ws.Range("E24").Value = "  -1.72%  "

# Row 25
$ws.Range("E25").Value = "  +0.31%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.87"
$ws.Range("E26").Value = "  -2.62%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.23"
$ws.Range("E27").Value = "  -1.38%  "

# Row 28
$ws.Range("E28").Value = "  -2.11%  "

# Row 29
$ws.Range("E29").Value = "  -1.99%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.74"
$ws.Range("E30").Value = "  -3.72%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.34"
$ws.Range("E31").Value = "  -2.38%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.28"
$ws.Range("E32").Value = "  -4.01%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0747"
$ws.Range("E34").Value = "  -2.22%  "

# Row 35
$ws.Range("E35").Value = "  -2.95%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.95"
$ws.Range("E36").Value = "  +1.84%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.88"
$ws.Range("E37").Value = "  -1.99%  "

# Row 38
$ws.Range("E38").Value = "  -0.74%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.78"
$ws.Range("E39").Value = "  -2.26%  "

# Row 40
$ws.Range("E40").Value = "  -2.26%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.91"
$ws.Range("E41").Value = "  -1.77%  "

# Row 42
$ws.Range("D42").Value = "2.000.48"
$ws.Range("E42").Value = "  +0.86%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.21"
$ws.Range("E43").Value = "  -9.97%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.56"
$ws.Range("E44").Value = "  -1.93%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0276"
$ws.Range("E45").Value = "  -3.06%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.92"
$ws.Range("E46").Value = "  -1.88%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.57"
$ws.Range("E47").Value = "  +4.27%  "

# Row 48
$ws.Range("D48").Value = "2.669.82"
$ws.Range("E48").Value = "  -1.39%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "94.96"
$ws.Range("E49").Value = "  -2.48%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.62"
$ws.Range("E50").Value = "  -0.66%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.87"
$ws.Range("E51").Value = "  -0.57%  "
